$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header: column B now holds "volume_price" instead of "material",
# and the old "price_per_cubic_meter" column (C) is removed.
$ws.Range("B1").Value = "volume_price"
$ws.Range("C1").ClearContents()

# Fill in the new data rows (product_id, volume_price)
$data = @(
    @(1, 12),
    @(8, 16),
    @(16, 13),
    @(20, 21),
    @(21, 8),
    @(23, 41),
    @(10, 24)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r++
}

# Widen columns B and C to match the edited layout
$ws.Columns.Item(2).ColumnWidth = 17.45
$ws.Columns.Item(3).ColumnWidth = 9.92

# Restore the active selection used when the file was last saved
$ws.Range("C11").Select()
